# "Weeks with and without input" — turn the old X/_ text grid into a
# colored weekly chronogram: a merged year header, a row of week labels,
# and a diagonal staircase of "input received" cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Sheet)
$ws.Name = "Sheet"

# Wipe the previous X/_ marker grid
$ws.Cells.Clear()

# Colors (VBA-style BGR-packed RGB() values)
$blue   = 192 * 65536 + 112 * 256 + 0      # RGB(0,112,192)
$orange = 0   * 65536 + 165 * 256 + 255    # RGB(255,165,0)
$white  = 255 * 65536 + 255 * 256 + 255    # RGB(255,255,255)

# Nine equally-sized week columns (B..J)
$ws.Range("B1:J1").ColumnWidth = 14.14

# Year header, merged across the week columns
$year = $ws.Range("B1:J1")
$year.Merge()
$ws.Range("B1").NumberFormat = "@"            # keep "2023" as text, not a number
$ws.Range("B1").Value = "2023"
$ws.Range("B1").Interior.Color = $blue
$ws.Range("B1").Font.Color = $white
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment = -4108     # xlCenter

# Week labels across row 2
$weeks = @("Week 1","Week 2","Week 3","Week 4","Week 5","Week 6","Week 7","Week 8","Week 9")
for ($i = 0; $i -lt $weeks.Length; $i++) {
    $cell = $ws.Cells.Item(2, 2 + $i)
    $cell.Value = $weeks[$i]
    $cell.Interior.Color = $blue
    $cell.Font.Color = $white
    $cell.HorizontalAlignment = -4108         # xlCenter
}

# Diagonal staircase of "input received" weeks (orange fill, no text)
$inputCells = @("B3", "C4", "D5", "E6", "F7", "G7", "H7", "I8")
foreach ($addr in $inputCells) {
    $ws.Range($addr).Interior.Color = $orange
}

Write-Output "chronogram updated"
